# Update existing values on rows 2-4, then delete rows 5-7 entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MatrixTest")

# Row 2 updates
$ws.Range("G2").Value = 8
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 8

# Row 3 updates
$ws.Range("G3").Value = 2
$ws.Range("J3").Value = 10
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 5.666666666666667

# Row 4 updates
$ws.Range("J4").Value = 4
$ws.Range("M4").Value = 8

# Remove rows 5, 6 and 7 (old arg1_2 block), shrinking the used range to A1:O4
$ws.Range("A5:A7").EntireRow.Delete()
